$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22 (shifts existing rows 22.. down to 23..)
$ws.Rows.Item(22).Insert()

# Populate the new row 22 with the epiworld-forecasts entry
$ws.Range("A22").Value = "epiworld-forecasts: Automatic Disease Forecasting with epiworldR"
$ws.Range("B22").Value = "epiworld-forecasts uses epiworldR, GitHub Actions, and Docker to generate disease forecasts that update automatically. While we provide an example forecast of COVID-19 case counts in Utah, this tool is an open-source, template repository that can easily be adapted to generate forecasts for different diseases."
$ws.Range("C22").Value = "Andrew Pulsipher"
$ws.Range("D22").Value = "a.pulsipher@utah.edu"
$ws.Range("E22").Value = "Yes"
$ws.Range("G22").Value = "Published"
$ws.Range("H22").Value = "MIT"
$ws.Range("I22").Value = "R"
$ws.Range("J22").Value = "Forecasters"
$ws.Range("K22").Value = "TBD"
$ws.Range("L22").Value = "Automation pipeline tool"
$ws.Range("M22").Value = "Data sources, model definition, calibration steps"
$ws.Range("N22").Value = "https://github.com/EpiForeSITE/epiworld-forecasts, https://epiforesite.github.io/epiworld-forecasts/, https://github.com/UofUEpiBio/epiworldR/"
$ws.Range("O22").Value = "https://github.com/EpiForeSITE/epiworld-forecasts"
